$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
# Halved/updated hours for two existing volunteers
$ws.Range("E4").Value = 0.25
$ws.Range("E8").Value = 0.25

# Two more volunteers now have hours recorded (previously blank)
$ws.Range("E43").Value = 0.25
$ws.Range("E67").Value = 0.25

# The SUM formula in E68 recalculates automatically after these writes.

# --- Number format change --------------------------------------------------
# Column E used a custom "0.0" format; switch it to the built-in "0.00" format
# for the whole column so every cell (filled or still blank) stays consistent.
$ws.Range("E1:E68").NumberFormat = "0.00"

# --- View / selection state -------------------------------------------------
# Scroll back to the top of the sheet and leave the selection on I27.
$ws.Range("I27").Select() | Out-Null
